$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that needs to be
# updated from 45172 (2023-09-03) to 45175 (2023-09-06) for every data
# row (rows 2 through 497).
$newDate = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0

$lastRow = 497
$range = $ws.Range("C2:C$lastRow")
$range.Value = $newDate
